# Update workbook/sheet title and the "through" date label from
# September 20 to September 21, then refresh the underlying daily
# carjacking counts (adds data collected for 2021-09-21 / the 2021-09-29
# data pull) for the affected neighborhood rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab
$ws.Name = "Through 2021-09-21"

# Update the column header text for the "current" month column
$ws.Range("B1").Value = "September 2021 (through September 21)"

# Garfield Park (row 2)
$ws.Range("B2").Value = 14
$ws.Range("AL2").Value = 4

# Austin (row 5)
$ws.Range("AC5").Value = 7
$ws.Range("AL5").Value = 4

# Auburn Gresham (row 7)
$ws.Range("T7").Value = 1

# West Town (row 10)
$ws.Range("B10").Value = 4

# Chatham (row 13)
$ws.Range("T13").Value = 3

# Grand Boulevard (row 18)
$ws.Range("AL18").Value = 4

# River North (row 21)
$ws.Range("AL21").Value = 1

# South Chicago (row 22)
$ws.Range("T22").Value = 2

# Ashburn (row 24)
$ws.Range("K24").Value = 1

# Avalon Park (row 27)
$ws.Range("K27").Value = 1

# Lincoln Park (row 30)
$ws.Range("B30").Value = 2

# Douglas (row 39)
$ws.Range("T39").Value = 1

# Ukrainian Village (row 42)
$ws.Range("T42").Value = 1

# Woodlawn (row 43)
$ws.Range("BD43").Value = 1

# Hermosa (row 54)
$ws.Range("AC54").Value = 1

# Dunning (row 67)
$ws.Range("B67").Value = 1

# East Side (row 68)
$ws.Range("AL68").Value = 1

# Pullman (row 91)
$ws.Range("B91").Value = 2

# Rogers Park (row 93)
$ws.Range("B93").Value = 3

# Rush & Division (row 94)
$ws.Range("AL94").Value = 1
